# Generate Report for Handback
# Refresh the handback-status timestamps that get stamped when the report
# is (re)generated: the "Latest HO Xliff Generate Date" on the Overview
# sheet, and the "Correspond Handoff Datetime" / "Correspond Handback
# DateTime" pair on each per-locale sheet for the c66060bc... row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 4 is the c66060bc-ffbe-40b7-80f7-70f5ec727e01.md file.
$wsOverview.Range("G4").Value = "2016-10-24 09:14:41"

# zh-cn sheet: row 4 is the c66060bc-ffbe-40b7-80f7-70f5ec727e01 handback.
$wsZhCn.Range("H4").Value = "2016-10-24 09:14:30"
$wsZhCn.Range("K4").Value = "2016-10-24 09:15:14"

# de-de sheet: row 4 is the c66060bc-ffbe-40b7-80f7-70f5ec727e01 handback.
$wsDeDe.Range("H4").Value = "2016-10-24 09:14:41"
$wsDeDe.Range("K4").Value = "2016-10-24 09:15:32"
